$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source published a new weekly batch of Melón / Tuna price rows for
# "Femacal de La Calera" ahead of the existing 2021-12-22 (serial 44552)
# block. Insert 3 new rows at 378 so the old rows 378-384 shift down to
# 381-387 unchanged, then populate the 3 new rows with the new data.
$ws.Rows("378:380").Insert()

# Row 378: Melón / Tuna / Extra / Provincia de Quillota (2022-02-03)
$ws.Range("A378").Value = 3
$ws.Range("B378").Value = 'Femacal de La Calera'
$ws.Range("C378").Value = 'Coquimbo'
$ws.Range("D378").Value = 44595
$ws.Range("E378").Value = 5
$ws.Range("F378").Value = 100112027
$ws.Range("G378").Value = 'Melón'
$ws.Range("H378").Value = 'Tuna'
$ws.Range("I378").Value = 'Extra'
$ws.Range("J378").Value = 350
$ws.Range("K378").Value = 1000
$ws.Range("L378").Value = 1000
$ws.Range("M378").Value = 1000
$ws.Range("N378").Value = '$/unidad'
$ws.Range("O378").Value = 'Provincia de Quillota'
$ws.Range("P378").Value = 1000
$ws.Range("Q378").Value = 1
$ws.Range("R378").Value = 'Hortaliza'

# Row 379: Melón / Tuna / Primera / Provincia de Quillota (2022-02-03)
$ws.Range("A379").Value = 3
$ws.Range("B379").Value = 'Femacal de La Calera'
$ws.Range("C379").Value = 'Coquimbo'
$ws.Range("D379").Value = 44595
$ws.Range("E379").Value = 5
$ws.Range("F379").Value = 100112027
$ws.Range("G379").Value = 'Melón'
$ws.Range("H379").Value = 'Tuna'
$ws.Range("I379").Value = 'Primera'
$ws.Range("J379").Value = 350
$ws.Range("K379").Value = 700
$ws.Range("L379").Value = 700
$ws.Range("M379").Value = 700
$ws.Range("N379").Value = '$/unidad'
$ws.Range("O379").Value = 'Provincia de Quillota'
$ws.Range("P379").Value = 700
$ws.Range("Q379").Value = 1
$ws.Range("R379").Value = 'Hortaliza'

# Row 380: Melón / Tuna / Segunda / Provincia de Quillota (2022-02-03)
$ws.Range("A380").Value = 3
$ws.Range("B380").Value = 'Femacal de La Calera'
$ws.Range("C380").Value = 'Coquimbo'
$ws.Range("D380").Value = 44595
$ws.Range("E380").Value = 5
$ws.Range("F380").Value = 100112027
$ws.Range("G380").Value = 'Melón'
$ws.Range("H380").Value = 'Tuna'
$ws.Range("I380").Value = 'Segunda'
$ws.Range("J380").Value = 330
$ws.Range("K380").Value = 500
$ws.Range("L380").Value = 500
$ws.Range("M380").Value = 500
$ws.Range("N380").Value = '$/unidad'
$ws.Range("O380").Value = 'Provincia de Quillota'
$ws.Range("P380").Value = 500
$ws.Range("Q380").Value = 1
$ws.Range("R380").Value = 'Hortaliza'
